# Update countries & provincias Spain
# - Refresh "Kenia" stats (row 102) with the latest case counts.
# - Re-sort the tied rows 198-200 (Belice / Santa Lucia / Nueva Caledonia),
#   209-211 (Montserrat / Seychelles / Groenlandia) and 214-216
#   (Sahara Occidental / San Bartolome / Bonaire...) to match updated order.
# - Bump the "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Timestamp footer.
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 15:35"

# Kenia refreshed totals.
Set-Row 102 @("Kenia", 1214, 22, 383, 780, 0, 1, 51)

# Belice / Santa Lucia / Nueva Caledonia now tied at 18 - re-ordered.
Set-Row 198 @("Belice", 18, 0, 16, 0, 0, 0, 2)
Set-Row 199 @("Santa Lucia", 18, 0, 18, 0, 0, 0, 0)
Set-Row 200 @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)

# Montserrat / Seychelles / Groenlandia now tied at 11 - re-ordered.
Set-Row 209 @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
Set-Row 210 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
Set-Row 211 @("Groenlandia", 11, 0, 11, 0, 0, 0, 0)

# Sahara Occidental / San Bartolome / Bonaire... now tied at 6 - re-ordered.
Set-Row 214 @("Sahara Occidental", 6, 0, 6, 0, 0, 0, 0)
Set-Row 215 @("San Bartolome", 6, 0, 6, 0, 0, 0, 0)
Set-Row 216 @("Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0)
